$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Every existing data row (2..494) had its "Förändrad" (column C) date
#    bumped from 45180 (2023-09-11) to 45181 (2023-09-12).
for ($r = 2; $r -le 494; $r++) {
    $ws.Cells.Item($r, 3).Value = 45181
}

# 2) Row 494 gained an explicit row height (ht="15" customHeight="1"),
#    matching the rest of the data rows.
$ws.Rows.Item(494).RowHeight = 15

# 3) Four brand-new rows (495-498) were appended with new cleaning
#    notifications for VÄSTERÅS / VÄSTMANLANDS LÄN.
$newRows = @(
    @{ Row = 495; A = "A 42257-2023"; G = 2.1 },
    @{ Row = 496; A = "A 42253-2023"; G = 3.1 },
    @{ Row = 497; A = "A 42251-2023"; G = 2 },
    @{ Row = 498; A = "A 42262-2023"; G = 3 }
)

foreach ($rowData in $newRows) {
    $r = $rowData.Row

    $ws.Cells.Item($r, 1).Value = $rowData.A          # A - Beteckning
    $ws.Cells.Item($r, 2).Value = 45180                # B - Datum
    $ws.Cells.Item($r, 2).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($r, 3).Value = 45181                # C - Förändrad
    $ws.Cells.Item($r, 3).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($r, 4).Value = "VÄSTMANLANDS LÄN"   # D - Län
    $ws.Cells.Item($r, 5).Value = "VÄSTERÅS"           # E - Kommun
    $ws.Cells.Item($r, 7).Value = $rowData.G           # G - Area (ha)
    $ws.Cells.Item($r, 8).Value = 0                    # H - Fridlysta
    $ws.Cells.Item($r, 9).Value = 0                    # I - Signalarter
    $ws.Cells.Item($r, 10).Value = 0                   # J - NT
    $ws.Cells.Item($r, 11).Value = 0                   # K - VU
    $ws.Cells.Item($r, 12).Value = 0                   # L - EN
    $ws.Cells.Item($r, 13).Value = 0                   # M - CR
    $ws.Cells.Item($r, 14).Value = 0                   # N - RE
    $ws.Cells.Item($r, 15).Value = 0                   # O - Rödlistade
    $ws.Cells.Item($r, 16).Value = 0                   # P - Hotade
    $ws.Cells.Item($r, 17).Value = 0                   # Q - Alla arter
    $ws.Cells.Item($r, 18).WrapText = $true            # R - Artnamn (empty, wrapped)
}

# Rows 495-497 keep the standard 15pt custom row height; row 498 (the
# last row) is left without an explicit height, exactly like row 494
# was before this edit.
$ws.Rows.Item(495).RowHeight = 15
$ws.Rows.Item(496).RowHeight = 15
$ws.Rows.Item(497).RowHeight = 15
